# ES006.xlsx — "View Dfn" sheet maintenance edit
#
# 1. Bump the Class Name (C10) and Revision timestamp (C14).
# 2. Mark listFg / pdfViewer field groups as editable (H35 / H36 = "yes").
# 3. Add formatting-only placeholder cells (blank cells carrying the
#    surrounding table's style) in column A for rows 20-26 and in column B
#    for rows 33-41 and 48-108, matching the rest of their tables.
# 4. Remove the obsolete "Plugin" field row (old row 73) from the Fields
#    table, which shifts every row below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("View Dfn")

# --- simple value updates ------------------------------------------------
$ws.Range("C10").Value = "es.views.es006.ES006"
$ws.Range("C14").Value = 20250819115100
$ws.Range("H35").Value = "yes"
$ws.Range("H36").Value = "yes"

# --- formatting-only blank cells (copy style from a sibling in the table) -

# Recordset table (rows 19-26): column A blank cells styled like A19
$ws.Range("A19").Copy()
$ws.Range("A20:A26").PasteSpecial(-4122)   # xlPasteFormats

# Field Groups / Fields tables: column B blank cells styled like B47
$ws.Range("B47").Copy()
$ws.Range("B33:B41").PasteSpecial(-4122)
$ws.Range("B48:B72").PasteSpecial(-4122)
$ws.Range("B74:B108").PasteSpecial(-4122)  # becomes B73:B107 after the row delete below

$excel.CutCopyMode = 0

# --- remove the obsolete Plugin field row --------------------------------
$ws.Rows.Item(73).Delete()
